# Auto commit: 2025-04-15 22:13:15
# Renames Statistics -> Daily_Statistics, inserts a new Annual_Statistics sheet
# (copied/derived from the daily one), updates the daily statistics values to the
# new dataset, recomputes the annual-statistics numbers, renames the VaR / ES
# labels, and refreshes the MVP_Stats numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the first sheet "Statistics" -> "Daily_Statistics"
# ---------------------------------------------------------------------------
$daily = $wb.Worksheets.Item(1)
$daily.Name = "Daily_Statistics"

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "Annual_Statistics" sheet right after Daily_Statistics,
#    cloning the layout/formatting of Daily_Statistics.
# ---------------------------------------------------------------------------
$annual = $wb.Worksheets.Add($null, $daily)
$annual.Name = "Annual_Statistics"
$daily.Range("A1:E11").Copy($annual.Range("A1:E11"))

# ---------------------------------------------------------------------------
# 3. Rename the VaR / Expected Shortfall labels on Daily_Statistics
#    (Annual_Statistics inherited the old labels via the copy, so fix both).
# ---------------------------------------------------------------------------
$daily.Range("A8").Value = "VaR 95%"
$daily.Range("A9").Value = "ES 95%"
$daily.Range("A10").Value = "VaR 99%"
$daily.Range("A11").Value = "ES 99%"

$annual.Range("A8").Value = "VaR 95%"
$annual.Range("A9").Value = "ES 95%"
$annual.Range("A10").Value = "VaR 99%"
$annual.Range("A11").Value = "ES 99%"

# ---------------------------------------------------------------------------
# 4. Refresh the Daily_Statistics numbers (B2:E11) with the updated dataset.
# ---------------------------------------------------------------------------
$daily.Range("B2").Value = 201.3
$daily.Range("C2").Value = 254.7
$daily.Range("D2").Value = 781.3
$daily.Range("E2").Value = 36.55

$daily.Range("B3").Value = 0.0004098013293508894
$daily.Range("C3").Value = 0.0006221043485396727
$daily.Range("D3").Value = 0.0002993732022570493
$daily.Range("E3").Value = -0.0001622489719508735

$daily.Range("B4").Value = 0.01382976152410843
$daily.Range("C4").Value = 0.01190981747785037
$daily.Range("D4").Value = 0.0127642524571061
$daily.Range("E4").Value = 0.03150980377232784

$daily.Range("B5").Value = 0.00019126230381371
$daily.Range("C5").Value = 0.0001418437523557101
$daily.Range("D5").Value = 0.0001629261407887392
$daily.Range("E5").Value = 0.0009928677337706061

$daily.Range("B6").Value = -0.1736236955280446
$daily.Range("C6").Value = -0.5310272909960784
$daily.Range("D6").Value = -0.4348566222261216
$daily.Range("E6").Value = -0.7210069636931663

$daily.Range("B7").Value = 1.355289324219858
$daily.Range("C7").Value = 5.926437800708015
$daily.Range("D7").Value = 6.508584345364187
$daily.Range("E7").Value = 6.518477693830253

$daily.Range("B8").Value = -0.02391909987079396
$daily.Range("C8").Value = -0.01736620301863545
$daily.Range("D8").Value = -0.02043358709577255
$daily.Range("E8").Value = -0.04868307573724512

$daily.Range("B9").Value = -0.03128387024261594
$daily.Range("C9").Value = -0.0281395131394166
$daily.Range("D9").Value = -0.03068661573676915
$daily.Range("E9").Value = -0.07146682517145642

$daily.Range("B10").Value = -0.03644107693407325
$daily.Range("C10").Value = -0.03352152295620385
$daily.Range("D10").Value = -0.03661309409350758
$daily.Range("E10").Value = -0.0752792562145345

$daily.Range("B11").Value = -0.0442008868805281
$daily.Range("C11").Value = -0.04985149434408001
$daily.Range("D11").Value = -0.04873649174613838
$daily.Range("E11").Value = -0.1242680872778185

# ---------------------------------------------------------------------------
# 5. Populate the Annual_Statistics sheet with the annualized numbers.
# ---------------------------------------------------------------------------
$annual.Range("B2").Value = 201.3
$annual.Range("C2").Value = 254.7
$annual.Range("D2").Value = 781.3
$annual.Range("E2").Value = 36.55

$annual.Range("B3").Value = 0.1024503323377223
$annual.Range("C3").Value = 0.1555260871349182
$annual.Range("D3").Value = 0.07484330056426232
$annual.Range("E3").Value = -0.04056224298771838

$annual.Range("B4").Value = 0.2186677295657215
$annual.Range("C4").Value = 0.1883107487344456
$annual.Range("D4").Value = 0.2018205519692798
$annual.Range("E4").Value = 0.4982137427276083

$annual.Range("B5").Value = 0.0478155759534275
$annual.Range("C5").Value = 0.03546093808892752
$annual.Range("D5").Value = 0.04073153519718479
$annual.Range("E5").Value = 0.2482169334426515

$annual.Range("B6").Value = -0.1736236955280446
$annual.Range("C6").Value = -0.5310272909960784
$annual.Range("D6").Value = -0.4348566222261216
$annual.Range("E6").Value = -0.7210069636931663

$annual.Range("B7").Value = 1.355289324219858
$annual.Range("C7").Value = 5.926437800708015
$annual.Range("D7").Value = 6.508584345364187
$annual.Range("E7").Value = 6.518477693830253

$annual.Range("B8").Value = -0.3781941758637406
$annual.Range("C8").Value = -0.2745837792388978
$annual.Range("D8").Value = -0.3230833799503321
$annual.Range("E8").Value = -0.7697470141608774

$annual.Range("B9").Value = -0.4946414199591536
$annual.Range("C9").Value = -0.4449247688439585
$annual.Range("D9").Value = -0.4851979970527825
$annual.Range("E9").Value = -1.129989723414279

$annual.Range("B10").Value = -0.5761840175054853
$annual.Range("C10").Value = -0.5300218158961246
$annual.Range("D10").Value = -0.5789038476077093
$annual.Range("E10").Value = -1.19026955100657

$annual.Range("B11").Value = -0.6988773857096181
$annual.Range("C11").Value = -0.7882213344514727
$annual.Range("D11").Value = -0.77059159541897
$annual.Range("E11").Value = -1.964850981352499

# ---------------------------------------------------------------------------
# 6. Update the MVP_Stats sheet (Portfolio Mean Return / Std Dev).
# ---------------------------------------------------------------------------
$mvpStats = $wb.Worksheets.Item("MVP_Stats")
$mvpStats.Range("B2").Value = 0.000222138666407758
$mvpStats.Range("C2").Value = 0.01108004749581097

# ---------------------------------------------------------------------------
# 7. Keep Daily_Statistics as the active/selected sheet (matches original file).
# ---------------------------------------------------------------------------
$daily.Activate()
